# Fruta / hortaliza, semanal
# Update the weekly price/date figures for the Haba subset (Agrícola del Norte S.A. de Arica).
# The edit reshuffles the Fecha (D), Volumen (J), Precio mínimo (K), Precio máximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) values across the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: D (date serial), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rows = @{
    2  = @{ D = 44243; J = 1200; K = 1200; L = 1300; M = 1250; P = 1250 }
    3  = @{ D = 44442; J = 1250; K = 850;  L = 900;  M = 875;  P = 875  }
    4  = @{ D = 44649; J = 600;  K = 900;  L = 1000; M = 950;  P = 950  }
    5  = @{ D = 44284; J = 1500; K = 800;  L = 850;  M = 825;  P = 825  }
    7  = @{ D = 44291; J = 1000; K = 1000; L = 1200; M = 1100; P = 1100 }
    8  = @{ D = 44229; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    9  = @{ D = 44341; J = 1300; K = 900;  L = 1000; M = 950;  P = 950  }
    10 = @{ D = 44175; J = 1600; K = 1000; L = 1200; M = 1100; P = 1100 }
    11 = @{ D = 44550; J = 1300; K = 1000; L = 1200; M = 1100; P = 1100 }
    12 = @{ D = 44476; J = 900;  K = 700;  L = 800;  M = 750;  P = 750  }
    13 = @{ D = 44453; J = 1000; K = 800;  L = 900;  M = 850;  P = 850  }
    14 = @{ D = 44607; J = 900;  K = 1300; L = 1400; M = 1350; P = 1350 }
    15 = @{ D = 44407; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 }
    16 = @{ D = 44449; J = 1300; K = 900;  L = 950;  M = 925;  P = 925  }
    17 = @{ D = 44638; J = 1000; K = 900;  L = 950;  M = 925;  P = 925  }
    18 = @{ D = 44455; J = 1100; K = 900;  L = 1000; M = 950;  P = 950  }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value2  = $vals.D   # Column D: Fecha
    $ws.Cells.Item($r, 10).Value2 = $vals.J   # Column J: Volumen
    $ws.Cells.Item($r, 11).Value2 = $vals.K   # Column K: Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $vals.L   # Column L: Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $vals.M   # Column M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value2 = $vals.P   # Column P: Precio $/Kg
}
